# Update on 20210731 画中人
#
# Replace English/Korean double quotation marks used for emphasis/quoted
# speech with single quotation marks in several Istina/May dialogue lines
# (en_US column C), and replace "의료팀" (medical team) with "의료부"
# (medical department) in a few Korean (ko_KR column D) lines.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = '[name="Istina"]  Is it recording? ''Test. Test.''
'
$ws.Range("C10").Value = '[name="Istina"]  Closure lent me this device for the purposes of... Well, the application says ''recording self-diagnosis and treatment of trauma using exposure therapy and related methods.''
'
$ws.Range("C17").Value = '[name="Istina"]  ''Why not try talking it through with yourself...?'' They made it sound so easy.
'
$ws.Range("C34").Value = '[name="Istina"]  I am called ''Istina''——This is, of course, a codename. My true name is Anna Morozova.
'
$ws.Range("C38").Value = '[name="Istina"]  Now I live at the headquarters of an organization called Rhodes Island, where I am classified as part of the ''Ursus Student Self-Government Group.''
'
$ws.Range("C43").Value = '[name="Istina"]  We currently serve as ''Operators'' for this Rhodes Island... I suppose it’s a kind of job.
'
$ws.Range("C239").Value = '[name="Istina"]  It’s been a while since I had a good talk about detective novels. And the plot in that one is indeed one of my favorites. ''Riveting,'' you said?
'
$ws.Range("C251").Value = '[name="Istina"]  You want to read too, Natalya? I never would have imagined. Every time I asked you before, you’d say the same thing: ''When I have time.''
'
$ws.Range("C265").Value = '[name="Istina"]  Also, what was that about a ''test,'' Natalya?
'
$ws.Range("C282").Value = '[name="May"]  And if we’re talking classics, we of course have to talk about ''that one,'' eh?
'
$ws.Range("C283").Value = '[name="Istina"]  ''That one?'' Do you mean...?
'
$ws.Range("C284").Value = '[name="May"]  ''That one!'' Yes!
'
$ws.Range("C286").Value = '[name="May"]  ''One truth prevails! The murderer is——You!''
'
$ws.Range("C383").Value = '[name="Istina"]  As a ''junior detective,'' keen powers of observation are a must, no?
'
$ws.Range("C412").Value = '[name="Istina"]  ''The good news,'' huh?
'
$ws.Range("D16").Value = '[name="이스티나"]  제 증상을 완화하고 치료하고 싶다고 이야기하니, 의료부의 오퍼레이터가 이런 제안을 해주었습니다……
'
$ws.Range("D262").Value = '[name="나탈리야"]  의료부 오퍼레이터들이 특히 더 그래. 얼마 전엔 테스트를 받으러 갔었는데, 메딕 몇 명이 박사의 영양 밸런스에 관해서 이야기하더라고…… 그때 뭐라고 했더라? 일단은 간식부터 금지해야 한다고 했던가?
'
$ws.Range("D394").Value = '[name="메이"]  드론을 쓰면 이렇게나 간단히 찾을 수 있는 것인가? 켈시 선생님부터 시작해서 의료부 사람들은 정말 대단한 것이다……
'
$ws.Range("D396").Value = '[name="나탈리야"]  미안, 안나. 의료부에 가봐야 할 거 같아.
'
